$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 13
$ws.Range("C3").Value = 6
$ws.Range("C4").Value = 11
$ws.Range("C5").Value = 19
$ws.Range("C7").Value = 21
$ws.Range("C8").Value = 18
$ws.Range("C9").Value = 13
$ws.Range("B10").Value = "<hind>"
$ws.Range("C10").Value = 15
$ws.Range("C12").Value = 12
$ws.Range("C13").Value = 14
$ws.Range("C14").Value = 11
$ws.Range("C15").Value = 16
$ws.Range("B16").Value = "<sulu>"
$ws.Range("C17").Value = 15
$ws.Range("B18").Value = "<in>"
